# Updates cryptos list values per upstream data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.335.42'
$ws.Range("E2").Value = '  +4.16%  '
$ws.Range("D3").Value = '1.580.40'
$ws.Range("E3").Value = '  +0.59%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -1.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.27'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.74%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.495'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.22%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.62'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +6.94%  '
$ws.Range("E9").Value = '  +0.74%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0600'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0886'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +2.03%  '
$ws.Range("D12").Value = '1.805.21'
$ws.Range("E12").Value = '  +0.60%  '
$ws.Range("D13").Value = '1.567.99'
$ws.Range("E13").Value = '  +4.54%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.77'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.63%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.528'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.34%  '
$ws.Range("D16").Value = '28.279.99'
$ws.Range("E16").Value = '  +3.98%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.83'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +2.37%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '231.91'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +7.24%  '
$ws.Range("D19").Value = '0.0₃0708'
$ws.Range("E19").Value = '  +0.77%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.46'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.41%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.998'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.13'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.47%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.35'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.95'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.17%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.94'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.33%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.27'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.19%  '
$ws.Range("E27").Value = '  -1.35%  '
$ws.Range("E28").Value = '  -0.37%  '
$ws.Range("E29").Value = '  -0.85%  '
$ws.Range("E30").Value = '  +0.18%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0473'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.27%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.23'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.49%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.15'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.93%  '
$ws.Range("D34").Value = '1.417.42'
$ws.Range("E34").Value = '  -2.39%  '
$ws.Range("E35").Value = '  -1.19%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.05'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -4.79%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.32'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.40%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.53'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +7.67%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.542'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.99%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.810'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.998'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.16%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.64'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -3.16%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.974'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.50%  '
$ws.Range("E45").Value = '  +5.18%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.18'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.76%  '
$ws.Range("D47").Value = '1.715.57'
$ws.Range("E47").Value = '  +0.58%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.25'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.48%  '
$ws.Range("E49").Value = '  +2.79%  '
$ws.Range("E50").Value = '  +1.02%  '
$ws.Range("B51").Value = 'BitcoinSV'
$ws.Range("C51").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '39.51'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +16.50%  '
